# Apply updated crypto price/volume figures to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) contain plain text values that can
# look like numbers/dates (e.g. "62.516.32", "0.506"). Force the affected
# cells to Text format first so Excel does not auto-convert them.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "62.516.32"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.91%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.181.01"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -3.71%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "586.93"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.94%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.24"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -4.99%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.178.04"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -3.76%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.506"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.71%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.141"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -5.77%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.24"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -4.32%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.453"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -3.91%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000234"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -5.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.26"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -4.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.707.50"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.65%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.92%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.187.57"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -3.50%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "62.577.02"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.92%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.51"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -5.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "455.55"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -5.22%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.89"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.77%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.701"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -4.18%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.60"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -5.34%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.50"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.43%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.23"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.11%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.89%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.08%  "
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = "NEARProtocol"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.82"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -6.51%  "
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.72"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -5.21%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.01"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -7.45%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.27"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -6.71%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.95%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -6.47%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.03"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -5.43%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.51%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "51.02"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -3.93%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0₃0696"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -6.97%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0384"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -4.26%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.70"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.31%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.97"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -4.90%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.821.62"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -7.78%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.12%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "384.68"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -9.53%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "36.18"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +3.51%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.248"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -6.54%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.06%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.12"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -3.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "124.87"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.61%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "25.34"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.63%  "
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -4.03%  "
